$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.761.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.530.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.79%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '624.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.93%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.614'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.526.74'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.77%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.199'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.03'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.85%  '

$ws.Range("E12").Value = '  -1.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000278'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.093.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.87%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '610.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.519.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.814.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.51%  '

$ws.Range("E20").Value = '  +0.99%  '

$ws.Range("E21").Value = '  +1.15%  '

$ws.Range("E22").Value = '  -0.95%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '98.52'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.78%  '

$ws.Range("E26").Value = '  -1.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("E28").Value = '  -3.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  -3.05%  '

$ws.Range("E31").Value = '  -2.47%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.28%  '

$ws.Range("E33").Value = '  -0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '638.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("E35").Value = '  -5.13%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.1000'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0472'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.63%  '

$ws.Range("E39").Value = '  -9.38%  '

$ws.Range("E40").Value = '  -0.85%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.144'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.74%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.351.88'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.03%  '

$ws.Range("E44").Value = '  +0.30%  '

$ws.Range("E45").Value = '  -1.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.312'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '31.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.97%  '

$ws.Range("E48").Value = '  -6.28%  '

$ws.Range("E49").Value = '  -0.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.78'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.80%  '

$ws.Range("E51").Value = '  -0.02%  '
